$wb = $excel.ActiveWorkbook

$sheetNames = @("Summary", "Pattern1-Pure Data")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Text-valued cells: force Text number format so Excel doesn't
    # auto-convert currency/percentage/date-looking strings into numbers.
    $textCells = @("D3", "E3", "F3", "G3", "I3", "J3", "K3", "L3", "O3")
    foreach ($cellRef in $textCells) {
        $ws.Range($cellRef).NumberFormat = "@"
    }

    $ws.Range("D3").Value = "¥1,001,002.00"
    $ws.Range("E3").Value = "¥+1,002.00"
    $ws.Range("F3").Value = "+0.10%"
    $ws.Range("G3").Value = "+28.71%"
    $ws.Range("H3").Value = 0
    $ws.Range("I3").Value = "0.00%"
    $ws.Range("J3").Value = "100.0%"
    $ws.Range("K3").Value = "0.1002%"
    $ws.Range("L3").Value = "0.0000%"
    $ws.Range("M3").Value = 2
    $ws.Range("O3").Value = "20251218"
}
